# Updates the cryptocurrency price/volume table to reflect the latest
# scrape (GitHub Actions symbol-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address -> new text value.
$updates = [ordered]@{
    'D2' = '307.97'
    'D3' = '40.94'
    'E3' = '1.98%'
    'D4' = '5.128'
    'E4' = '0.63%'
    'D5' = '0.07616'
    'E5' = '-1.09%'
    'D6' = '1.626'
    'E6' = '0.50%'
    'E7' = '0.04%'
    'D8' = '0.8997'
    'E8' = '2.51%'
    'D9' = '0.1086'
    'E9' = '9.48%'
    'D10' = '0.1761'
    'E10' = '1.77%'
    'D11' = '0.09200'
    'E11' = '2.64%'
    'D12' = '0.04205'
    'E12' = '-4.51%'
    'D13' = '0.1051'
    'E13' = '-0.49%'
    'D14' = '0.001254'
    'E14' = '-0.23%'
    'D15' = '0.005910'
    'E15' = '-0.31%'
    'B16' = 'HotbitToken'
    'C16' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D16' = '0.004095'
    'E16' = '0.75%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.352'
    'E17' = '-0.12%'
    'B18' = 'GateToken'
    'C18' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D18' = '4.253'
    'E18' = '-0.29%'
    'B19' = 'BitpandaEcosystemToken'
    'C19' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D19' = '0.3294'
    'E19' = '-0.20%'
    'B20' = 'MCDex'
    'C20' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D20' = '6.568'
    'E20' = '-6.07%'
    'B21' = 'ProBitToken'
    'C21' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D21' = '0.1364'
    'E21' = '1.92%'
    'B22' = 'ZBToken'
    'C22' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D22' = '0.2681'
    'E22' = '-14.53%'
    'B23' = 'CoinExToken'
    'C23' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D23' = '0.04088'
    'E23' = '-1.54%'
    'B24' = 'BitKan'
    'C24' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D24' = '0.001222'
    'E24' = '1.99%'
    'E25' = '6.60%'
    'E38' = '1.18%'
    'D39' = '0.05174'
    'E39' = '0.59%'
    'D40' = '0.007760'
    'E40' = '-2.17%'
    'D41' = '0.1299'
    'E41' = '-1.67%'
    'D42' = '0.006783'
    'E42' = '6.56%'
    'D43' = '0.001951'
    'E43' = '0.09%'
    'E44' = '-0.44%'
    'D45' = '0.3075'
    'E45' = '0.86%'
    'D46' = '0.00006936'
    'E46' = '6.40%'
    'E47' = '0.05%'
    'D48' = '0.03238'
    'E48' = '852.56%'
    'D49' = '0.004202'
    'E49' = '-39.97%'
    'D50' = '0.00002101'
    'E50' = '0.05%'
    'D51' = '0.0002001'
    'E51' = '0.05%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format first so values like "307.97" or "1.98%" are
    # written as plain text (matching the workbook's inlineStr cells)
    # rather than being auto-converted to numbers/percentages by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop the temporary text-format styling so the cell keeps using
    # the default (unstyled) formatting, same as before the edit.
    $cell.ClearFormats()
}
